# "Generate Report for Handback": the handback CI run refreshed the
# "Correspond Handoff/Handback Datetime" stamps (columns E/H) for the first
# data row of each language-report sheet (row 2), leaving the second data
# row (row 3) untouched.
$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# zh-cn report: row 2 (0e314636-...zh-cn.xlf) got a new handoff/handback pass
$wsZhCn.Range("E2").Value = "2016-03-22 08:46:37"
$wsZhCn.Range("H2").Value = "2016-03-22 08:47:00"

# de-de report: row 2 (0e314636-...de-de.xlf) got a new handoff/handback pass
$wsDeDe.Range("E2").Value = "2016-03-22 08:46:41"
$wsDeDe.Range("H2").Value = "2016-03-22 08:47:08"
